$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be pre-formatted as
# Text so Excel keeps them as strings (matching the source inlineStr cells)
# instead of silently converting them to floating point numbers.
$textCells = $excel.Union($ws.Range("D5"), $ws.Range("D6"), $ws.Range("D9"), $ws.Range("D10"), $ws.Range("D13"), $ws.Range("D14"), $ws.Range("D18"), $ws.Range("D20"), $ws.Range("D23"), $ws.Range("D25"), $ws.Range("D26"), $ws.Range("D27"), $ws.Range("D29"), $ws.Range("D32"), $ws.Range("D34"), $ws.Range("D35"), $ws.Range("D36"), $ws.Range("D37"), $ws.Range("D38"), $ws.Range("D41"), $ws.Range("D43"), $ws.Range("D44"), $ws.Range("D45"), $ws.Range("D46"), $ws.Range("D47"), $ws.Range("D48"), $ws.Range("D49"), $ws.Range("D50"), $ws.Range("D51"))
foreach ($c in $textCells) {
    $c.NumberFormat = "@"
}

$ws.Range('D2').Value = '66.126.92'
$ws.Range('E2').Value = '  +1.23%  '
$ws.Range('D3').Value = '2.691.78'
$ws.Range('E3').Value = '  +1.63%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '609.43'
$ws.Range('E5').Value = '  +0.72%  '
$ws.Range('D6').Value = '158.88'
$ws.Range('E6').Value = '  +0.83%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -0.56%  '
$ws.Range('D9').Value = '0.127'
$ws.Range('E9').Value = '  +6.94%  '
$ws.Range('D10').Value = '6.06'
$ws.Range('E10').Value = '  +3.17%  '
$ws.Range('E11').Value = '  +0.82%  '
$ws.Range('E12').Value = '  +1.56%  '
$ws.Range('D13').Value = '0.0000211'
$ws.Range('E13').Value = '  +17.06%  '
$ws.Range('D14').Value = '30.28'
$ws.Range('E14').Value = '  +3.08%  '
$ws.Range('D15').Value = '3.177.29'
$ws.Range('E15').Value = '  +1.57%  '
$ws.Range('D16').Value = '65.983.08'
$ws.Range('E16').Value = '  +1.22%  '
$ws.Range('D17').Value = '2.680.56'
$ws.Range('E17').Value = '  +1.99%  '
$ws.Range('D18').Value = '12.77'
$ws.Range('E18').Value = '  +1.14%  '
$ws.Range('E19').Value = '  +1.35%  '
$ws.Range('D20').Value = '363.17'
$ws.Range('E20').Value = '  +2.18%  '
$ws.Range('E21').Value = '  +2.42%  '
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('D23').Value = '70.19'
$ws.Range('E23').Value = '  +2.88%  '
$ws.Range('E24').Value = '  +1.69%  '
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').Value = '0.0000108'
$ws.Range('E25').Value = '  +12.57%  '
$ws.Range('B26').Value = 'SuiNetwork'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D26').Value = '1.68'
$ws.Range('E26').Value = '  -2.18%  '
$ws.Range('D27').Value = '1.70'
$ws.Range('E27').Value = '  +1.53%  '
$ws.Range('E28').Value = '  +4.52%  '
$ws.Range('D29').Value = '8.27'
$ws.Range('E29').Value = '  -0.37%  '
$ws.Range('E30').Value = '  +4.50%  '
$ws.Range('E31').Value = '  +0.30%  '
$ws.Range('D32').Value = '534.41'
$ws.Range('E32').Value = '  -1.65%  '
$ws.Range('E33').Value = '  -0.76%  '
$ws.Range('D34').Value = '6.66'
$ws.Range('E34').Value = '  +2.76%  '
$ws.Range('D35').Value = '5.48'
$ws.Range('E35').Value = '  -5.68%  '
$ws.Range('D36').Value = '0.435'
$ws.Range('E36').Value = '  +1.04%  '
$ws.Range('D37').Value = '20.85'
$ws.Range('E37').Value = '  +3.01%  '
$ws.Range('D38').Value = '163.20'
$ws.Range('E38').Value = '  -1.24%  '
$ws.Range('E39').Value = '  -2.57%  '
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('D41').Value = '170.90'
$ws.Range('E41').Value = '  +1.13%  '
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('D43').Value = '42.95'
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('D44').Value = '4.19'
$ws.Range('E44').Value = '  +1.21%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').Value = '0.0616'
$ws.Range('E45').Value = '  +0.89%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').Value = '2.32'
$ws.Range('E46').Value = '  +1.88%  '
$ws.Range('D47').Value = '23.49'
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('D48').Value = '0.661'
$ws.Range('E48').Value = '  +1.12%  '
$ws.Range('D49').Value = '0.0266'
$ws.Range('E49').Value = '  +5.27%  '
$ws.Range('D50').Value = '20.48'
$ws.Range('E50').Value = '  +4.55%  '
$ws.Range('D51').Value = '0.0988'
$ws.Range('E51').Value = '  +0.24%  '
